# "Order Confirmation" sheet (xl/worksheets/sheet9.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Confirmation")

# Row 7: TS_ORD_01 / FR-ORD-01 description changes
$ws.Range("C7").Value = "Verify order display on confrimation page"

# Row 11: TS_ORD_05 / FR-ORD-05 (new)
$ws.Range("A11").Value = "TS_ORD_05"
$ws.Range("B11").Value = "FR-ORD-05"

# Row 9: TS_ORD_03 / FR-ORD-03 description changes
$ws.Range("C9").Value = "Verify order display details"

# Row 10: TS_ORD_04 / FR-ORD-04 description changes
$ws.Range("C10").Value = "Verify order display payment status"

# Row 11: description
$ws.Range("C11").Value = "Verify order confrimation message"

# Row 12: TS_ORD_06 (new)
$ws.Range("A12").Value = "TS_ORD_06"

# Row 12: description
$ws.Range("C12").Value = "Verify order in MyOrder"

# Row 13: TS_ORD_07 (new)
$ws.Range("A13").Value = "TS_ORD_07"

# Row 12: FR-ORD-06 (new)
$ws.Range("B12").Value = "FR-ORD-06"

# Row 13: FR-ORD-07 (new)
$ws.Range("B13").Value = "FR-ORD-07"

# Row 13: description
$ws.Range("C13").Value = "Verfiy Email/Phone order confrimation"

# --- Activate "Order Confirmation" tab and move the selection ---
$ws.Activate()
$ws.Range("C16").Select()
